# Swap the contents of column B ("category") and column C ("option")
# on every worksheet, per commit: "Update order of parameters in elic_cat
# (topic > option > category > ...)"

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $lastRow = $used.Row + $used.Rows.Count - 1
    $firstRow = 1

    $colB = $ws.Range($ws.Cells.Item($firstRow, 2), $ws.Cells.Item($lastRow, 2))
    $colC = $ws.Range($ws.Cells.Item($firstRow, 3), $ws.Cells.Item($lastRow, 3))

    $valuesB = $colB.Value2
    $valuesC = $colC.Value2

    $colB.Value2 = $valuesC
    $colC.Value2 = $valuesB
}
